# Update time_taken timestamps in the "data" sheet (column F, rows 2-17)
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

$dataSheet.Range("F2").Value = "2021-10-05 14:34:34.817932"
$dataSheet.Range("F3").Value = "2021-10-05 14:34:34.817941"
$dataSheet.Range("F4").Value = "2021-10-05 14:34:34.817944"
$dataSheet.Range("F5").Value = "2021-10-05 14:34:34.817947"
$dataSheet.Range("F6").Value = "2021-10-05 14:34:34.817951"
$dataSheet.Range("F7").Value = "2021-10-05 14:34:34.817954"
$dataSheet.Range("F8").Value = "2021-10-05 14:34:34.817956"
$dataSheet.Range("F9").Value = "2021-10-05 14:34:34.817959"
$dataSheet.Range("F10").Value = "2021-10-05 14:34:34.817962"
$dataSheet.Range("F11").Value = "2021-10-05 14:34:34.817965"
$dataSheet.Range("F12").Value = "2021-10-05 14:34:34.817968"
$dataSheet.Range("F13").Value = "2021-10-05 14:34:34.817970"
$dataSheet.Range("F14").Value = "2021-10-05 14:34:34.817973"
$dataSheet.Range("F15").Value = "2021-10-05 14:34:34.817975"
$dataSheet.Range("F16").Value = "2021-10-05 14:34:34.817978"
$dataSheet.Range("F17").Value = "2021-10-05 14:34:34.817980"

# Add the new "metadata" sheet right after "data"
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$ws.Name = "metadata"

# Header row
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Data row
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Maturity-onset Diabetes of the Young"
$ws.Range("C2").Value = 301
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.10"
$ws.Range("E2").Value = "2021-03-01T05:15:38.475360Z"
$ws.Range("F2").Value = "2021-10-05 14:34:34.814520"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/301/?format=json"

# Apply the same header style used on the "data" sheet's header row / index column
$dataSheet.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats
$dataSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

# Match the page margins used throughout the rest of the workbook
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

Write-Output "metadata sheet added and timestamps updated"
